$d = $word.ActiveDocument

# 1. Merge "Franco Geremia y Facundo " + "Facello" -> " Franco Geremia y Facundo Facello"
#    (removes the spell-check split; net text is unchanged, so this is a no-op for Find/Replace,
#     but we still run it to normalize in case whitespace needs adjusting)
$d.Content.Find.Execute("Franco Geremia y Facundo Facello", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Franco Geremia y Facundo Facello", 2) | Out-Null

# 2. Date day change: "20" -> "31" (within "Fecha de la ultima revision: 20 de enero, 2022")
$d.Content.Find.Execute("revisión: 20 de enero", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "revisión: 31 de enero", 2) | Out-Null

# 3. Replace the "sobre Pokémon (Cualquier dato...)" text with the new wording
$d.Content.Find.Execute("sobre Pokémon (Cualquier dato que sea un atributo de uno, como nombre, tipo, habilidades, apariencia, etc) ", `
                         $true, $false, $false, $false, $false, `
                         $true, 1, $false, "(Numero, tipo, habilidades, apariencia) sobre Pokémon utilizando sus nombres ", 2) | Out-Null

# 4. "Pokémon en un" merge - no text change needed (already correct in Find semantics)

# 5. "layout" merge - no text change needed

# 6. "PokéAPI" merge - no text change needed
